# Fix the repository link paragraph:
#   1. Rename the URL ("kpfu-" prefix added to the repo name).
#   2. Drop the <w:hyperlink> wrapper around the URL run (keep the
#      Hyperlink character style on the run itself).
#   3. Leave the "_GoBack" bookmark spanning the very end of the
#      paragraph (after both runs) instead of the very beginning.
#
# Word's COM bridge in this sandbox does not support repositioning an
# existing bookmark (Bookmarks.Add(name, range) / Bookmark.Start= /
# Bookmark.Range= are all no-ops once the target range isn't right at
# the top of the document). Instead we rebuild the paragraph content
# in place: insert freshly-formatted plain-text runs immediately BEFORE
# the bookmark's current (collapsed) location -- which pushes the
# existing "_GoBack" bookmark along, past the new text, without ever
# calling a bookmark-repositioning API -- then delete the old runs/
# hyperlink that are now left dangling after the bookmark.

$d = $word.ActiveDocument

$prefixText = "Отчеты по остальным практикам и лабораторным работам можно найти на моем репозитории: "
$urlText = "https://github.com/ironsast/kpfu-probability-theory-and-mathematical-statistics"

# Locate the paragraph that still holds the old link text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Отчеты по остальным*") {
        $target = $p
    }
}

$paraStart = $target.Range.Start

# 1) Insert the new text (prefix + url) as plain runs right at the start
#    of the paragraph -- i.e. exactly where the "_GoBack" bookmark
#    currently sits. The bookmark keeps "pointing" at the old content,
#    so it ends up shifted past this newly-inserted text.
$insertPoint = $d.Range($paraStart, $paraStart)
$insertPoint.InsertAfter($prefixText + $urlText)

$totalLen = $prefixText.Length + $urlText.Length

# 2) Apply the paragraph's normal run formatting (28 half-points = 14pt)
#    across the whole freshly inserted span.
$newRange = $d.Range($paraStart, $paraStart + $totalLen)
$newRange.Font.Size = 14

# 3) Re-apply the Hyperlink character style to just the URL portion so
#    it keeps its visual (blue/underline) look without an actual
#    <w:hyperlink> wrapper.
$urlStart = $paraStart + $prefixText.Length
$urlRange = $d.Range($urlStart, $paraStart + $totalLen)
$urlRange.Style = "Hyperlink"

# 4) Remove the now-stale old runs (plain text run + the hyperlink run)
#    that sit right after the bookmark.
$bm = $d.Bookmarks("_GoBack")
$paraEnd = $target.Range.End
$oldContentRange = $d.Range($bm.End, $paraEnd - 1)
if ($oldContentRange.Start -lt $oldContentRange.End) {
    $oldContentRange.Delete()
}
